$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("F1").Value = "United States"

# Updated values for columns B through F, rows 2-10
$data = @(
    @(0.686274509803922, 0.689716312056738, 0.804147465437788, 0.667359667359667, 0.550151975683891),
    @(0.803455723542117, 0.823333333333333, 0.850678733031674, 0.798850574712644, 0.677653902084343),
    @(0.741721854304636, 0.79020979020979,  0.775229357798165, 0.709090909090909, 0.592648539778449),
    @(0.528436018957346, 0.435452793834297, 0.601965601965602, 0.587301587301587, 0.452247191011236),
    @(0.666666666666667, 0.695364238410596, 0.760526315789474, 0.6953125,         0.569230769230769),
    @(0.595693779904306, 0.718095238095238, 0.786885245901639, 0.844036697247706, 0.627764127764128),
    @(0.792941176470588, 0.794117647058824, 0.807228915662651, 0.816120906801008, 0.630372492836676),
    @(0.891304347826087, 0.847790507364976, 0.906040268456376, 0.87378640776699,  0.6211714132187),
    @(0.57906976744186,  0.608247422680412, 0.618556701030928, 0.667359667359667, 0.465809768637532)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rowData[0]
    $ws.Cells.Item($row, 3).Value = $rowData[1]
    $ws.Cells.Item($row, 4).Value = $rowData[2]
    $ws.Cells.Item($row, 5).Value = $rowData[3]
    $ws.Cells.Item($row, 6).Value = $rowData[4]
}
